$wb = $excel.ActiveWorkbook

# --- Sheet "Data" ---
$wsData = $wb.Worksheets.Item("Data")

$wsData.Range("A2").Value = 3024
$wsData.Range("E2").Value = 46200608024
$wsData.Range("X2").Value = "DN4127460130024"

$wsData.Range("A3").Value = 3025
$wsData.Range("E3").Value = 46200608025
$wsData.Range("X3").Value = "DN4127460130025"

# --- Sheet "Check" ---
$wsCheck = $wb.Worksheets.Item("Check")

$wsCheck.Range("A2").Value = 3024
$wsCheck.Range("C2").Value = "DN4127460130024"

$wsCheck.Range("A3").Value = 3025
$wsCheck.Range("C3").Value = "DN4127460130025"
